$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from
# 45180 (2023-09-11) to 45181 (2023-09-12) for every data row (2..252).
for ($r = 2; $r -le 252; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = 45181
}
